# Fall 2021 workshop update: explicitly stamp every paragraph (body +
# paragraph styles) with PageBreakBefore = False, matching the source
# edit's <w:pageBreakBefore w:val="0"/> addition throughout the package.

$d = $word.ActiveDocument

# Body paragraphs: every <w:p><w:pPr> in the document gets an explicit
# pageBreakBefore element.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs($i).Format.PageBreakBefore = 0
}

# Paragraph styles (Heading 1-6, Title, Subtitle) also get the explicit
# pageBreakBefore element added to their stored pPr.
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $d.Styles($name).ParagraphFormat.PageBreakBefore = 0
}

Write-Host "Applied pageBreakBefore=0 to $($d.Paragraphs.Count) paragraphs and $($styleNames.Count) styles"
